$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The values in columns D, M, N, O, P, S for rows 3-7 have been
# cyclically re-ordered. New row 3 gets old row 7's values, new row 4
# gets old row 5's values, new row 5 gets old row 6's values, new row 6
# gets old row 4's values, and new row 7 gets old row 3's values.

$data = @{
    3 = @{ D = 44229; M = 55; N = 11000; O = 12000; P = 11364; S = 812 }
    4 = @{ D = 44216; M = 55; N = 11000; O = 12000; P = 11545; S = 825 }
    5 = @{ D = 44172; M = 90; N = 8500;  O = 9000;  P = 8806;  S = 629 }
    6 = @{ D = 44210; M = 70; N = 10000; O = 11000; P = 10357; S = 740 }
    7 = @{ D = 44232; M = 60; N = 11000; O = 12000; P = 11583; S = 827 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
